# Auto-generated Word COM-interop script to apply the diff
# (updates the worksheet date and the 100 answer cells in the table)
$d = $word.ActiveDocument

# 1. Update the date line at the top of the document.
#    wdReplaceOne (1) so only the single matching occurrence is touched.
$d.Content.Find.Execute("2025-04-15 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-22 Tuesday", 1) | Out-Null

# 2. Update each answer cell in the table (row, col) -> new text.
$t = $d.Tables.Item(1)
$replacements = @(
    @{Row=1; Col=1; Old="10+74=84"; New="12+54=66"},
    @{Row=1; Col=2; Old="34+65=99"; New="15+82=97"},
    @{Row=1; Col=3; Old="91-64=27"; New="24+45=69"},
    @{Row=1; Col=4; Old="18+77=95"; New="54+5=59"},
    @{Row=1; Col=5; Old="85-25=60"; New="16+20=36"},
    @{Row=2; Col=1; Old="13+15=28"; New="56+13=69"},
    @{Row=2; Col=2; Old="75+22=97"; New="44+8=52"},
    @{Row=2; Col=3; Old="93-48=45"; New="95-37=58"},
    @{Row=2; Col=4; Old="5+54=59"; New="65-6=59"},
    @{Row=2; Col=5; Old="72-35=37"; New="78-26=52"},
    @{Row=3; Col=1; Old="53-41=12"; New="5+28=33"},
    @{Row=3; Col=2; Old="23+23=46"; New="37+26=63"},
    @{Row=3; Col=3; Old="12+10=22"; New="52+45=97"},
    @{Row=3; Col=4; Old="96-83=13"; New="86-14=72"},
    @{Row=3; Col=5; Old="66-34=32"; New="96-32=64"},
    @{Row=4; Col=1; Old="56+2=58"; New="25+18=43"},
    @{Row=4; Col=2; Old="55-20=35"; New="3+64=67"},
    @{Row=4; Col=3; Old="3+23=26"; New="85-1=84"},
    @{Row=4; Col=4; Old="96-6=90"; New="6+62=68"},
    @{Row=4; Col=5; Old="40+34=74"; New="91-34=57"},
    @{Row=5; Col=1; Old="21+9=30"; New="1+70=71"},
    @{Row=5; Col=2; Old="23+43=66"; New="95-89=6"},
    @{Row=5; Col=3; Old="83-38=45"; New="47+21=68"},
    @{Row=5; Col=4; Old="70-51=19"; New="83-17=66"},
    @{Row=5; Col=5; Old="93-28=65"; New="77-27=50"},
    @{Row=6; Col=1; Old="59-33=26"; New="24-2=22"},
    @{Row=6; Col=2; Old="86-3=83"; New="32-20=12"},
    @{Row=6; Col=3; Old="26+8=34"; New="37+56=93"},
    @{Row=6; Col=4; Old="3-1=2"; New="38+37=75"},
    @{Row=6; Col=5; Old="14+34=48"; New="37-31=6"},
    @{Row=7; Col=1; Old="5+46=51"; New="68-36=32"},
    @{Row=7; Col=2; Old="84-55=29"; New="77-28=49"},
    @{Row=7; Col=3; Old="34+50=84"; New="73+20=93"},
    @{Row=7; Col=4; Old="31+26=57"; New="97-42=55"},
    @{Row=7; Col=5; Old="20+64=84"; New="33-5=28"},
    @{Row=8; Col=1; Old="98-65=33"; New="7+87=94"},
    @{Row=8; Col=2; Old="24+68=92"; New="39+42=81"},
    @{Row=8; Col=3; Old="18-5=13"; New="55-9=46"},
    @{Row=8; Col=4; Old="95-18=77"; New="49+29=78"},
    @{Row=8; Col=5; Old="86-24=62"; New="30-4=26"},
    @{Row=9; Col=1; Old="51+28=79"; New="48-6=42"},
    @{Row=9; Col=2; Old="25+74=99"; New="8+66=74"},
    @{Row=9; Col=3; Old="75-36=39"; New="6+61=67"},
    @{Row=9; Col=4; Old="66-37=29"; New="64-52=12"},
    @{Row=9; Col=5; Old="24+1=25"; New="60-47=13"},
    @{Row=10; Col=1; Old="49+33=82"; New="87-19=68"},
    @{Row=10; Col=2; Old="23+39=62"; New="80-64=16"},
    @{Row=10; Col=3; Old="21+60=81"; New="78-9=69"},
    @{Row=10; Col=4; Old="12+37=49"; New="93-56=37"},
    @{Row=10; Col=5; Old="30+43=73"; New="7+1=8"},
    @{Row=11; Col=1; Old="72+4=76"; New="88-35=53"},
    @{Row=11; Col=2; Old="69-60=9"; New="59-31=28"},
    @{Row=11; Col=3; Old="15+42=57"; New="98-16=82"},
    @{Row=11; Col=4; Old="43-3=40"; New="43+44=87"},
    @{Row=11; Col=5; Old="66+13=79"; New="88-5=83"},
    @{Row=12; Col=1; Old="50-21=29"; New="98-24=74"},
    @{Row=12; Col=2; Old="17+47=64"; New="63-60=3"},
    @{Row=12; Col=3; Old="64-50=14"; New="58+18=76"},
    @{Row=12; Col=4; Old="30+47=77"; New="69-6=63"},
    @{Row=12; Col=5; Old="84-54=30"; New="95-16=79"},
    @{Row=13; Col=1; Old="49-24=25"; New="49+4=53"},
    @{Row=13; Col=2; Old="78-6=72"; New="4+39=43"},
    @{Row=13; Col=3; Old="20+16=36"; New="81-17=64"},
    @{Row=13; Col=4; Old="35+44=79"; New="11+69=80"},
    @{Row=13; Col=5; Old="83-72=11"; New="70-41=29"},
    @{Row=14; Col=1; Old="1+16=17"; New="69-8=61"},
    @{Row=14; Col=2; Old="21-4=17"; New="12+72=84"},
    @{Row=14; Col=3; Old="28+38=66"; New="69-63=6"},
    @{Row=14; Col=4; Old="96-40=56"; New="71-54=17"},
    @{Row=14; Col=5; Old="62+10=72"; New="88+2=90"},
    @{Row=15; Col=1; Old="64-11=53"; New="42+37=79"},
    @{Row=15; Col=2; Old="25+16=41"; New="81+18=99"},
    @{Row=15; Col=3; Old="74-34=40"; New="11+64=75"},
    @{Row=15; Col=4; Old="43-36=7"; New="47+19=66"},
    @{Row=15; Col=5; Old="82-73=9"; New="49-25=24"},
    @{Row=16; Col=1; Old="82-74=8"; New="57-32=25"},
    @{Row=16; Col=2; Old="78-50=28"; New="33+20=53"},
    @{Row=16; Col=3; Old="80+11=91"; New="41+2=43"},
    @{Row=16; Col=4; Old="34+6=40"; New="45-20=25"},
    @{Row=16; Col=5; Old="72-50=22"; New="31+36=67"},
    @{Row=17; Col=1; Old="34+28=62"; New="56+32=88"},
    @{Row=17; Col=2; Old="91-88=3"; New="50+34=84"},
    @{Row=17; Col=3; Old="21+30=51"; New="87-28=59"},
    @{Row=17; Col=4; Old="95-5=90"; New="78-3=75"},
    @{Row=17; Col=5; Old="67+26=93"; New="11+35=46"},
    @{Row=18; Col=1; Old="15+15=30"; New="82+3=85"},
    @{Row=18; Col=2; Old="85+2=87"; New="52-32=20"},
    @{Row=18; Col=3; Old="98-34=64"; New="89-17=72"},
    @{Row=18; Col=4; Old="32-18=14"; New="62-18=44"},
    @{Row=18; Col=5; Old="48+15=63"; New="23+37=60"},
    @{Row=19; Col=1; Old="22+56=78"; New="29+7=36"},
    @{Row=19; Col=2; Old="25+16=41"; New="72-1=71"},
    @{Row=19; Col=3; Old="30+65=95"; New="5+9=14"},
    @{Row=19; Col=4; Old="89-70=19"; New="91-42=49"},
    @{Row=19; Col=5; Old="59-41=18"; New="74+20=94"},
    @{Row=20; Col=1; Old="56-5=51"; New="89-42=47"},
    @{Row=20; Col=2; Old="47+44=91"; New="11+69=80"},
    @{Row=20; Col=3; Old="70-43=27"; New="75-45=30"},
    @{Row=20; Col=4; Old="82-27=55"; New="8+39=47"},
    @{Row=20; Col=5; Old="69-66=3"; New="48+23=71"}
)

foreach ($item in $replacements) {
    $cell = $t.Cell($item.Row, $item.Col)
    $cellRange = $cell.Range
    # wdReplaceOne (1): the search text can repeat elsewhere in the table
    # (e.g. two cells both originally read "25+16=41"), so scope the Find to
    # this cell's Range and only replace the single occurrence found there.
    $cellRange.Find.Execute($item.Old, $true, $false, $false, $false, $false, $true, 1, $false, $item.New, 1) | Out-Null
}

Write-Host "Done applying replacements."
